# Apply "sound and mushroom height" edits to the active sheet (Sheet1).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 12 (id=8): volume 1 -> 0.5, isStereo 1 -> 0, falloffDistance 3000 -> 600
$ws.Range("E12").Value = 0.5
$ws.Range("G12").Value = 0
$ws.Range("I12").Value = 600

# Row 14 (id=10): soundGuid 202244 -> 201870, isStereo 1 -> 0, falloffDistance 3000 -> 600
$ws.Range("B14").Value = 201870
$ws.Range("G14").Value = 0
$ws.Range("I14").Value = 600

# Row 15 (id=11): isStereo 1 -> 0, falloffDistance 3000 -> 600
$ws.Range("G15").Value = 0
$ws.Range("I15").Value = 600

# Row 16 (id=12): isStereo 1 -> 0, falloffDistance 3000 -> 600
$ws.Range("G16").Value = 0
$ws.Range("I16").Value = 600

# Row 17 (id=13): isStereo 1 -> 0, falloffDistance 3000 -> 600
$ws.Range("G17").Value = 0
$ws.Range("I17").Value = 600

# Move the active selection to E17 (was E26) as recorded in the sheet view.
$ws.Activate()
$ws.Range("E17").Select()
